$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new ward columns (Ramani Huria wards "Mburahati" and "Ubungo") ---
# Before the edit, column layout (row 1) is:
#   ... K=Manzese, L=Mchikichini ... S=Temeke, T=Vingunguti
# Insert "Mburahati" right after Manzese (before Mchikichini) at column L.
$ws.Columns("L").Insert()
# After that insert, Temeke is now column T and Vingunguti is column U.
# Insert "Ubungo" right after Temeke (before Vingunguti) at column U.
$ws.Columns("U").Insert()

# --- Fill in the header row for the two new columns ---
$ws.Range("L1").Value = "Mburahati"
$ws.Range("U1").Value = "Ubungo"

# --- Fill in the new ward data for each data row ---
$ws.Range("L2").Value = 27
$ws.Range("U2").Value = 96

$ws.Range("L3").Value = 155
$ws.Range("U3").Value = 126

$ws.Range("L4").Value = 5
$ws.Range("U4").Value = 15

$ws.Range("L5").Value = 22
$ws.Range("U5").Value = 81

# --- Remove the stray selection row (row 7, used to hold the active-cell marker) ---
$ws.Rows(7).Delete()

# --- Reset the active selection back to the top-left cell now that A7 is gone ---
$ws.Range("A1").Select() | Out-Null
